$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 302
$startSerial = 44376
$count = 27

for ($i = 0; $i -lt $count; $i++) {
    $r = $startRow + $i
    $serial = $startSerial + $i

    # Duplicate the format of the row above (same style used throughout column A)
    $srcRow = $ws.Range("A" + ($r - 1) + ":D" + ($r - 1))
    $dstRow = $ws.Range("A" + $r + ":D" + $r)
    $srcRow.Copy($dstRow)

    $ws.Cells.Item($r, 1).Value = $serial
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
}
